$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2561.4614
$ws.Range("I28").Value = 1006.375
$ws.Range("K28").Value = 1006.375
$ws.Range("M28").Value = -521.375
$ws.Range("H64").Value = 28575084
$ws.Range("I64").Value = 3786.6562
$ws.Range("J64").Value = 333335580
$ws.Range("K64").Value = 3786.6562
$ws.Range("L64").Value = 333335580
$ws.Range("M64").Value = -3538.6562
$ws.Range("N64").Value = -333336076
$ws.Range("H67").Value = 28575084
$ws.Range("I67").Value = 3786.6562
$ws.Range("J67").Value = 333335580
$ws.Range("K67").Value = 3786.6562
$ws.Range("L67").Value = 333335580
$ws.Range("M67").Value = -2928.6562
$ws.Range("N67").Value = -333337296
$ws.Range("H88").Value = 10001049
$ws.Range("I88").Value = 25000698
$ws.Range("J88").Value = 1283.1666
$ws.Range("K88").Value = 25000698
$ws.Range("L88").Value = 1283.1666
$ws.Range("M88").Value = -25000292
$ws.Range("N88").Value = -2095.1666
$ws.Range("H91").Value = 10001049
$ws.Range("I91").Value = 25000698
$ws.Range("J91").Value = 1283.1666
$ws.Range("K91").Value = 25000698
$ws.Range("L91").Value = 1283.1666
$ws.Range("M91").Value = -24999294
$ws.Range("N91").Value = -4091.1666
$ws.Range("H111").Value = 1163.3334
$ws.Range("I111").Value = 229
$ws.Range("K111").Value = 687
$ws.Range("M111").Value = 2380
$ws.Range("H118").Value = 673.38464
$ws.Range("I118").Value = 632.63635
$ws.Range("K118").Value = 1897.90905
$ws.Range("M118").Value = -240.90905
$ws.Range("H137").Value = 2508040
$ws.Range("I137").Value = 10001140
$ws.Range("K137").Value = 30003420
$ws.Range("M137").Value = -30000870
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7355.227
$ws.Range("I32").Value = 6937.6313
$ws.Range("K32").Value = 6937.6313
$ws.Range("M32").Value = -6650.6313
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 120000
$ws.Range("J122").Value = 120000
$ws.Range("L122").Value = 120000
$ws.Range("N122").Value = -129800
$ws.Range("H134").Value = 2690.1914
$ws.Range("I134").Value = 2394.3125
$ws.Range("K134").Value = 7182.9375
$ws.Range("M134").Value = -4647.9375
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4706.2383
$ws.Range("I31").Value = 3327.8
$ws.Range("J31").Value = 8152.3335
$ws.Range("K31").Value = 3327.8
$ws.Range("L31").Value = 8152.3335
$ws.Range("M31").Value = -3032.8
$ws.Range("N31").Value = -8742.333500000001
$ws.Range("H34").Value = 4706.2383
$ws.Range("I34").Value = 3327.8
$ws.Range("J34").Value = 8152.3335
$ws.Range("K34").Value = 3327.8
$ws.Range("L34").Value = 8152.3335
$ws.Range("M34").Value = -3125.8
$ws.Range("N34").Value = -8556.333500000001
$ws.Range("H58").Value = 3846.4614
$ws.Range("I58").Value = 3584.3333
$ws.Range("J58").Value = 4071.1428
$ws.Range("K58").Value = 3584.3333
$ws.Range("L58").Value = 4071.1428
$ws.Range("M58").Value = -3381.3333
$ws.Range("N58").Value = -4477.1428
$ws.Range("H86").Value = 13954.625
$ws.Range("J86").Value = 9463.666999999999
$ws.Range("L86").Value = 9463.666999999999
$ws.Range("N86").Value = -11709.667
$ws.Range("H89").Value = 13954.625
$ws.Range("J89").Value = 9463.666999999999
$ws.Range("L89").Value = 47318.335
$ws.Range("N89").Value = -58550.335
$ws.Range("H99").Value = 26666.334
$ws.Range("I99").Value = 26666.334
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 26666.334
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -25168.334
$ws.Range("H107").Value = 1054.05
$ws.Range("I107").Value = 737.8
$ws.Range("K107").Value = 737.8
$ws.Range("M107").Value = 1182.2
$ws.Range("H126").Value = 26666.334
$ws.Range("I126").Value = 26666.334
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 79999.00199999999
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -77529.00199999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 5481.0586
$ws.Range("I132").Value = 4863.1816
$ws.Range("J132").Value = 6613.8335
$ws.Range("K132").Value = 14589.5448
$ws.Range("L132").Value = 19841.5005
$ws.Range("M132").Value = -12059.5448
$ws.Range("N132").Value = -24901.5005
$ws.Range("H134").Value = 4431.273
$ws.Range("I134").Value = 5124.1665
$ws.Range("K134").Value = 15372.4995
$ws.Range("M134").Value = -12837.4995
$ws.Range("H135").Value = 75000
$ws.Range("J135").Value = 75000
$ws.Range("L135").Value = 75000
$ws.Range("N135").Value = -85140
$ws.Range("H136").Value = 3846.4614
$ws.Range("I136").Value = 3584.3333
$ws.Range("J136").Value = 4071.1428
$ws.Range("K136").Value = 10752.9999
$ws.Range("L136").Value = 12213.4284
$ws.Range("M136").Value = -8202.999899999999
$ws.Range("N136").Value = -17313.4284
$ws.Range("H2").Value = 1528.4
$ws.Range("J2").Value = 3018.9
$ws.Range("L2").Value = 18113.4
$ws.Range("N2").Value = -18339.4
$ws.Range("H26").Value = 892.95
$ws.Range("J26").Value = 870.9091
$ws.Range("L26").Value = 2612.7273
$ws.Range("N26").Value = -3188.7273
$ws.Range("H68").Value = 9101330
$ws.Range("I68").Value = 9000
$ws.Range("K68").Value = 27000
$ws.Range("M68").Value = -26189
$ws.Range("H71").Value = 9101330
$ws.Range("I71").Value = 9000
$ws.Range("K71").Value = 81000
$ws.Range("M71").Value = -76944
$ws.Range("H75").Value = 1966
$ws.Range("J75").Value = 1966
$ws.Range("L75").Value = 5898
$ws.Range("N75").Value = -7894
$ws.Range("H78").Value = 1966
$ws.Range("J78").Value = 1966
$ws.Range("L78").Value = 17694
$ws.Range("N78").Value = -27678
$ws.Range("H107").Value = 427.8889
$ws.Range("I107").Value = 390
$ws.Range("K107").Value = 1170
$ws.Range("M107").Value = 750
$ws.Range("H121").Value = 150540.75
$ws.Range("I121").Value = 50150
$ws.Range("J121").Value = 184004.33
$ws.Range("K121").Value = 150450
$ws.Range("L121").Value = 552012.99
$ws.Range("M121").Value = -149140
$ws.Range("N121").Value = -554632.99
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H70").Value = 85282.28
$ws.Range("I70").Value = 157943.31
$ws.Range("J70").Value = 6566.1665
$ws.Range("K70").Value = 157943.31
$ws.Range("L70").Value = 6566.1665
$ws.Range("M70").Value = -157673.31
$ws.Range("N70").Value = -7106.1665
$ws.Range("H73").Value = 85282.28
$ws.Range("I73").Value = 157943.31
$ws.Range("J73").Value = 6566.1665
$ws.Range("K73").Value = 157943.31
$ws.Range("L73").Value = 6566.1665
$ws.Range("M73").Value = -157007.31
$ws.Range("N73").Value = -8438.166499999999
$ws.Range("H102").Value = 9622.166999999999
$ws.Range("J102").Value = 13788.875
$ws.Range("L102").Value = 13788.875
$ws.Range("N102").Value = -17032.875
$ws.Range("H126").Value = 8685.286
$ws.Range("I126").Value = 2100
$ws.Range("J126").Value = 17465.666
$ws.Range("K126").Value = 6300
$ws.Range("L126").Value = 52396.99800000001
$ws.Range("M126").Value = -3830
$ws.Range("N126").Value = -57336.99800000001
$ws.Range("H132").Value = 3600.6
$ws.Range("I132").Value = 2644.9
$ws.Range("K132").Value = 7934.700000000001
$ws.Range("M132").Value = -5404.700000000001
$ws.Range("H136").Value = 14531.4
$ws.Range("J136").Value = 14531.4
$ws.Range("L136").Value = 43594.2
$ws.Range("N136").Value = -48694.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 545.86365
$ws.Range("I55").Value = 260.69232
$ws.Range("J55").Value = 957.7778
$ws.Range("K55").Value = 260.69232
$ws.Range("L55").Value = 957.7778
$ws.Range("M55").Value = -87.69232
$ws.Range("N55").Value = -1303.7778
$ws.Range("H93").Value = 603
$ws.Range("I93").Value = 603
$ws.Range("K93").Value = 603
$ws.Range("M93").Value = 645
$ws.Range("H132").Value = 5052.25
$ws.Range("J132").Value = 6699.5713
$ws.Range("L132").Value = 20098.7139
$ws.Range("N132").Value = -25158.7139
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 53980
$ws.Range("J43").Value = 53930
$ws.Range("L43").Value = 53930
$ws.Range("N43").Value = -54228
$ws.Range("H81").Value = 1955.9445
$ws.Range("I81").Value = 2153.1
$ws.Range("K81").Value = 4306.2
$ws.Range("M81").Value = -3245.2
$ws.Range("H84").Value = 1955.9445
$ws.Range("I84").Value = 2153.1
$ws.Range("K84").Value = 21531
$ws.Range("M84").Value = -16227
$ws.Range("H122").Value = 7813830.5
$ws.Range("I122").Value = 1296.6957
$ws.Range("K122").Value = 3890.0871
$ws.Range("M122").Value = -1440.0871
$ws.Range("H132").Value = 4885
$ws.Range("I132").Value = 4994
$ws.Range("J132").Value = 4757.8335
$ws.Range("K132").Value = 14982
$ws.Range("L132").Value = 14273.5005
$ws.Range("M132").Value = -12452
$ws.Range("N132").Value = -19333.5005
$ws.Range("H137").Value = 78999.10000000001
$ws.Range("J137").Value = 78999.10000000001
$ws.Range("L137").Value = 78999.10000000001
$ws.Range("N137").Value = -89199.10000000001
